$p = $ppt.ActivePresentation
$s2 = $p.Slides.Item(2)

# 1) Duplicate slide 2 (the original 12-picture slide) -> becomes slide 3
$dup = $s2.Duplicate()

# 2) Reposition / relabel the 12 pictures already on slide 2 into the new 3x4 grid
$sh = $s2.Shapes.Item(1)
$sh.Name = "Picture 4"
$sh.Left = 406.15655517578125
$sh.Top = 103.745361328125
$sh.Width = 293.8099365234375
$sh.Height = 246.66709899902344

$sh = $s2.Shapes.Item(2)
$sh.Name = "Picture 6"
$sh.Left = 102.495361328125
$sh.Top = 103.74567413330078
$sh.Width = 293.8097839355469
$sh.Height = 246.66685485839844

$sh = $s2.Shapes.Item(3)
$sh.Name = "Picture 8"
$sh.Left = 1013.4793090820312
$sh.Top = 103.74543762207031
$sh.Width = 293.8096923828125
$sh.Height = 246.66685485839844

$sh = $s2.Shapes.Item(4)
$sh.Name = "Picture 10"
$sh.Left = 709.8180541992188
$sh.Top = 103.74551391601562
$sh.Width = 293.80963134765625
$sh.Height = 246.66685485839844

$sh = $s2.Shapes.Item(5)
$sh.Name = "Picture 21"
$sh.Left = 102.0114974975586
$sh.Top = 609.18212890625
$sh.Width = 292.7861633300781
$sh.Height = 139.93685913085938

$sh = $s2.Shapes.Item(6)
$sh.Name = "Picture 23"
$sh.Left = 407.18048095703125
$sh.Top = 609.18212890625
$sh.Width = 292.7860107421875
$sh.Height = 139.93685913085938

$sh = $s2.Shapes.Item(7)
$sh.Name = "Picture 25"
$sh.Left = 708.1112670898438
$sh.Top = 609.18212890625
$sh.Width = 292.78607177734375
$sh.Height = 139.93685913085938

$sh = $s2.Shapes.Item(8)
$sh.Name = "Picture 35"
$sh.Left = 1011.772705078125
$sh.Top = 609.18212890625
$sh.Width = 295.5163269042969
$sh.Height = 141.2417449951172

$sh = $s2.Shapes.Item(9)
$sh.Name = "Picture 20"
$sh.Left = 102.495361328125
$sh.Top = 366.56097412109375
$sh.Width = 292.3022155761719
$sh.Height = 226.47268676757812

$sh = $s2.Shapes.Item(10)
$sh.Name = "Picture 22"
$sh.Left = 404.64898681640625
$sh.Top = 364.2125244140625
$sh.Width = 293.8099365234375
$sh.Height = 228.82102966308594

$sh = $s2.Shapes.Item(11)
$sh.Name = "Picture 24"
$sh.Left = 708.1112670898438
$sh.Top = 366.5608825683594
$sh.Width = 293.8097839355469
$sh.Height = 226.47276306152344

$sh = $s2.Shapes.Item(12)
$sh.Name = "Picture 26"
$sh.Left = 1011.772705078125
$sh.Top = 366.5608825683594
$sh.Width = 295.5163269042969
$sh.Height = 226.47276306152344

# 3) Add the column/row caption textboxes
$tb = $s2.Shapes.AddTextbox(1, 205.13449096679688, 71.84481048583984, 88.53126525878906, 31.504724502563477)
$tb.TextFrame.WordWrap = 0
$tb.TextFrame.AutoSize = 1
$tb.TextFrame.TextRange.Text = "1 second"
$tb.TextFrame.TextRange.Font.Bold = -1
$tb.TextFrame.TextRange.Font.Size = 20

$tb = $s2.Shapes.AddTextbox(1, 504.7568664550781, 71.51433563232422, 96.609375, 31.504724502563477)
$tb.TextFrame.WordWrap = 0
$tb.TextFrame.AutoSize = 1
$tb.TextFrame.TextRange.Text = "3 seconds"
$tb.TextFrame.TextRange.Font.Bold = -1
$tb.TextFrame.TextRange.Font.Size = 20

$tb = $s2.Shapes.AddTextbox(1, 806.1995849609375, 72.24055480957031, 96.609375, 31.504724502563477)
$tb.TextFrame.WordWrap = 0
$tb.TextFrame.AutoSize = 1
$tb.TextFrame.TextRange.Text = "5 seconds"
$tb.TextFrame.TextRange.Font.Bold = -1
$tb.TextFrame.TextRange.Font.Size = 20

$tb = $s2.Shapes.AddTextbox(1, 1111.22607421875, 72.24055480957031, 96.609375, 31.504724502563477)
$tb.TextFrame.WordWrap = 0
$tb.TextFrame.AutoSize = 1
$tb.TextFrame.TextRange.Text = "7 seconds"
$tb.TextFrame.TextRange.Font.Bold = -1
$tb.TextFrame.TextRange.Font.Size = 20

$tb = $s2.Shapes.AddTextbox(1, 55.50952911376953, 211.3264617919922, 37.13417434692383, 31.504724502563477)
$tb.TextFrame.WordWrap = 0
$tb.TextFrame.AutoSize = 1
$tb.TextFrame.TextRange.Text = "(a)"
$tb.TextFrame.TextRange.Font.Bold = -1
$tb.TextFrame.TextRange.Font.Size = 20

$tb = $s2.Shapes.AddTextbox(1, 55.50952911376953, 462.8706359863281, 38.01763916015625, 31.504724502563477)
$tb.TextFrame.WordWrap = 0
$tb.TextFrame.AutoSize = 1
$tb.TextFrame.TextRange.Text = "(b)"
$tb.TextFrame.TextRange.Font.Bold = -1
$tb.TextFrame.TextRange.Font.Size = 20

$tb = $s2.Shapes.AddTextbox(1, 55.50952911376953, 663.398193359375, 35.619529724121094, 31.504724502563477)
$tb.TextFrame.WordWrap = 0
$tb.TextFrame.AutoSize = 1
$tb.TextFrame.TextRange.Text = "(c)"
$tb.TextFrame.TextRange.Font.Bold = -1
$tb.TextFrame.TextRange.Font.Size = 20

# 4) On the duplicated slide (now slide 3), drop all but the first picture
$s3 = $p.Slides.Item(3)
for ($i = $s3.Shapes.Count; $i -ge 2; $i--) {
    $s3.Shapes.Item($i).Delete() | Out-Null
}

